# Applies the commit's spreadsheet changes:
#  - Column header A1: "Gen" -> "MaxFES"
#  - Column A (rows 2-14): generation counters -> MaxFES fractions
#  - Drop the last "Run 50" data column (old AZ) and the old "Mean" column (old BA);
#    the new last column (AZ) becomes the recomputed "Mean" of the remaining 50 runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column A (MaxFES fractions), rows 2..14
$newA = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)

# New recomputed Mean values (of the 50 remaining runs), rows 2..14
$newMean = @(13.91551517, 12.72058333, 9.99682917, 5.52410968, 3.82288988, 2.70497334, 2.05601397, 1.71193678, 1.47112475, 1.32816028, 1.26264813, 1.23313558, 1.22557702)

# 1. Remove the entire last column (BA): this was the old "Mean" column.
#    After this, the former "Run 50" column (AZ) becomes the new last column.
$ws.Range("BA1:BA14").Delete()

# 2. Update header and column A data values
$ws.Range("A1").Value = "MaxFES"
for ($i = 0; $i -lt $newA.Count; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $newA[$i]
}

# 3. Replace the former "Run 50" column (AZ) with the recomputed "Mean" column
$ws.Range("AZ1").Value = "Mean"
for ($i = 0; $i -lt $newMean.Count; $i++) {
    $row = $i + 2
    $ws.Range("AZ$row").Value = $newMean[$i]
}
